# LoginData.xlsx edit: trim the LogIn sheet down to the first data row,
# point the AF login row at the new env2 host/db, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hyperlink on C3 (old af_pl_env1 URL) is going away along with that URL,
# so drop the hyperlink object before the cell value gets overwritten.
$ws.Hyperlinks.Delete()

# Rows 4-7 ("LogIn Alis 3/4/5/6") are removed entirely, leaving just the
# header row plus the two remaining test cases (TFL, AF).
$ws.Rows("4:7").Delete()

# Row 3 ("LogIn Alis AF") now points at the new alf-app01 env2 host/db.
$ws.Range("C3").Value = "http://alis-alf-app01:8082/af_pl_env2/alis#alis"
$ws.Range("F3").Value = "af_7000_michael"

# Final selection in the saved file is F3.
[void]$ws.Range("F3").Select()
